# Weekly price-sheet update: a new daily record is inserted above the
# current row 23, pushing the existing rows 23-32 down to 24-33 (the
# historical rows keep their original data, just shifted down by one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 23; rows 23:32 shift down to 24:33.
$ws.Rows.Item(23).Insert()

# Fill the new row 23 with the latest weekly record for
# "Terminal Hortofrutícola Agro Chillán" / Haba.
$ws.Cells.Item(23, 1).Value = 7
$ws.Cells.Item(23, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(23, 3).Value = "Ñuble"
$ws.Cells.Item(23, 4).Value = 44529
$ws.Cells.Item(23, 5).Value = 16
$ws.Cells.Item(23, 6).Value = 100112026
$ws.Cells.Item(23, 7).Value = "Haba"
$ws.Cells.Item(23, 8).Value = "Sin especificar"
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 100
$ws.Cells.Item(23, 11).Value = 6000
$ws.Cells.Item(23, 12).Value = 7000
$ws.Cells.Item(23, 13).Value = 6500
$ws.Cells.Item(23, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(23, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(23, 16).Value = 260
$ws.Cells.Item(23, 17).Value = 25
$ws.Cells.Item(23, 18).Value = "Hortaliza"
